# edit.ps1
#
# Repairs the broken "{{RG}" merge field (missing closing brace) so it
# reads "{{RG}}", and tidies up the "São Carlos, {{DATA}}" run split
# (merges "São Carlos" and ", " into a single run) on slide 1 of the
# certificate template.
#
# Target shape: slide 1 / "Rectangle 5" (the big paragraph of body text
# that contains the {{NOME}}, {{RG}}, {{CPF}}, {{DATA}} merge fields).

$p  = $ppt.ActivePresentation
$s  = $p.Slides.Item(1)

# Locate the "Rectangle 5" shape defensively (rather than hard-coding an
# index) by scanning the shapes on slide 1 for the one with a text frame
# whose text contains the "{{RG}" marker.
$sh = $null
for ($i = 1; $i -le $s.Shapes.Count; $i++) {
    $cand = $s.Shapes.Item($i)
    if ($cand.HasTextFrame) {
        if ($cand.TextFrame.TextRange.Text.IndexOf("{{RG}") -ge 0) {
            $sh = $cand
            break
        }
    }
}

$tr = $sh.TextFrame.TextRange

# ---------------------------------------------------------------------
# Fix 1: "{{RG} " (a single trailing brace - a typo) -> "{{RG}} " split
# across three runs: "{{", "RG", "}} ".
# ---------------------------------------------------------------------
$full    = $tr.Text
$rgIdx0  = $full.IndexOf("{{RG}")          # 0-based offset of "{"
$rgStart = $rgIdx0 + 1                     # 1-based COM Characters() start

# Re-type "RG" in place (splits "{{RG" into "{{" + "RG").
$midRun = $tr.Characters($rgStart + 2, 2)
$midRun.Text = "RG"

# Re-type the trailing "} " as "}} " (splits off the final run and adds
# the missing closing brace).
$tailRun = $tr.Characters($rgStart + 4, 2)
$tailRun.Text = "}} "

# ---------------------------------------------------------------------
# Fix 2: merge the "São Carlos" run and the ", " run into a single run
# "São Carlos, " (leaves the following "{{DATA}}" run untouched).
# ---------------------------------------------------------------------
$full2      = $tr.Text
$carlosIdx0 = $full2.IndexOf("Carlos")
$saoStart   = $carlosIdx0 - 3               # back up over "S","ã","o"," "... 

$mergedRun = $tr.Characters($saoStart, 12)  # "São Carlos, " (12 characters)
$mergedRun.Text = "São Carlos, "
